$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.349.50"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").Value = "1.668.03"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5350"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.72%  "
$ws.Range("E7").Value = "  +0.82%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2664"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06403"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07855"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("E12").Value = "  +1.58%  "
$ws.Range("D13").Value = "1.674.44"
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("D14").Value = "1.896.29"
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5541"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").Value = "0.0₅8200"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").Value = "26.367.98"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.690"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.68%  "
$ws.Range("E21").Value = "  +1.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.051"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1234"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.219"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.503"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05866"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.288"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.649"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.290"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.85%  "
$ws.Range("E34").Value = "  +1.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9729"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.831"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.420"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5848"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01602"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8646"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.59%  "
$ws.Range("D41").Value = "1.064.02"
$ws.Range("E41").Value = "  +3.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.843"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.75%  "
$ws.Range("E43").Value = "  +0.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("D45").Value = "1.807.61"
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.62%  "
$ws.Range("E47").Value = "  -5.38%  "
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4386"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.988"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.99%  "
$ws.Range("E51").Value = "  +0.47%  "
